# Add the MLFQ limitations bullet (and trailing blank paragraph) to the
# speaker notes of slide 7 ("Algorithm Implementations"), matching the
# commit "Add arguments for mlfq".
#
# Before:
#   1. Turns out basic algorithms only need a sorted queue on process data
#   2. The more advanced algorithms need their own data structure
#
# After:
#   1. Turns out basic algorithms only need a sorted queue on process data
#   2. The more advanced algorithms need their own data structure
#   3. MLFQ is unverified.  The implementation is limited to FIFO queues and
#      Round Robin algorithm.  The quantum for each queue is doubled.  There
#      is no starvation prevention of promoting processes.
#   4. (blank trailing paragraph)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(7)
$notes = $slide.NotesPage

# Shape 1 is the slide-image placeholder, shape 2 is the notes body text box.
$notesBody = $notes.Shapes.Item(2)
$tr = $notesBody.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
$lines = @()
for ($i = 1; $i -le $paraCount; $i++) {
    $lines += $tr.Paragraphs($i, 1).Text
}

$lines += "MLFQ is unverified.  The implementation is limited to FIFO queues and Round Robin algorithm.  The quantum for each queue is doubled.  There is no starvation prevention of promoting processes."
$lines += ""

$tr.Text = ($lines -join "`n")
